{"js": "// Update the year in the astromap link from 2019 to 2022.\n// \"Jenik Hollan, CzechGlobe (http://.../GaNight/2019/).\" -> \"...GaNight/2022/).\"\nconst oldText =\n  \"Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2019/).\";\nconst newText =\n  \"Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the whole sentence (all 4 differently-formatted runs) with a\n  // single plain run containing the updated year.\n  const target = results.items[0];\n  target.clear();\n  await context.sync();\n\n  target.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n} else {\n  // Fallback: only the year substring could be found (formatting/content\n  // already partly edited) - just swap the year in place.\n  const yearResults = body.search(\"GaNight/2019\", { matchCase: true });\n  yearResults.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < yearResults.items.length; i++) {\n    yearResults.items[i].insertText(\"GaNight/2022\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the year in the astromap link from 2019 to 2022.\n# \"Jenik Hollan, CzechGlobe (http://.../GaNight/2019/).\" -> \"...GaNight/2022/).\"\n$d = $word.ActiveDocument\n\n$oldText = \"Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2019/).\"\n$newText = \"Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = $oldText\n$found = $rng.Find.Execute()\n\nif ($found) {\n    # Replace the whole sentence (all differently-formatted runs) with a\n    # single plain run containing the updated year.\n    $rng.Text = \"\"\n    $rng.InsertAfter($newText)\n} else {\n    # Fallback: only the year substring could be found (e.g. formatting or\n    # content was already partly edited) - just swap the year in place.\n    $rng2 = $d.Content\n    $rng2.Find.ClearFormatting()\n    $rng2.Find.Text = \"GaNight/2019\"\n    $rng2.Find.Replacement.ClearFormatting()\n    $rng2.Find.Replacement.Text = \"GaNight/2022\"\n    $rng2.Find.Execute([ref]\"GaNight/2019\", $false, $false, $false, $false, $false, $true, 1, $false, \"GaNight/2022\", 2)\n}\n"}
